# "Generate Report for Handback" -- records the handback (target/xlf) files
# and timestamps for the two localized files, for both the zh-cn and de-de
# languages, and flips every row's Status from "Ready for handoff" to
# "Handed back: in sync with en-US". Also widens a few columns that now
# need to show longer handback file names.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item(1)
$zhcn     = $wb.Worksheets.Item(2)
$dede     = $wb.Worksheets.Item(3)

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 1. Status column: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (every row, every sheet, currently shows the same shared text)
# ---------------------------------------------------------------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2. Column widths -- columns that now hold long handback file names / ids
#    need to be widened. (ColumnWidth is in "characters"; the engine snaps
#    the stored xlsx width to the nearest 1/6 character, same as Excel's
#    own pixel-grid snapping, so we pick the ColumnWidth whose stored
#    width lands closest to the target.)
# ---------------------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668   # E
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668   # F

$zhcn.Columns.Item(3).ColumnWidth  = 29.166666666666668      # C  Status
$zhcn.Columns.Item(9).ColumnWidth  = 39.166666666666664      # I  Latest Target File
$zhcn.Columns.Item(10).ColumnWidth = 39.166666666666664      # J  Latest Handback File

$dede.Columns.Item(3).ColumnWidth  = 29.166666666666668      # C  Status
$dede.Columns.Item(9).ColumnWidth  = 39.166666666666664      # I  Latest Target File
$dede.Columns.Item(10).ColumnWidth = 39.166666666666664      # J  Latest Handback File

# ---------------------------------------------------------------------
# 3. zh-cn sheet: record target (.md) file, handback (.xlf) file and the
#    handback datetime for both handed-off documents.
# ---------------------------------------------------------------------
$zhcnHandbackTime = "2016-08-29 14:33:19"

$zhcn.Range("I2").Value = "a0e94d00-7c96-44aa-b06a-c48ee498b831.md"
$zhcn.Range("I2").Style = "Hyperlink"
$zhcn.Range("J2").Value = "a0e94d00-7c96-44aa-b06a-c48ee498b831.b1e5499c5bead7fe5de11bda301bce57ce7596d1.zh-cn.xlf"
$zhcn.Range("K2").Value = $zhcnHandbackTime

$zhcn.Range("I3").Value = "ee11e692-79d8-43d4-89c4-7343d3a4f709.md"
$zhcn.Range("I3").Style = "Hyperlink"
$zhcn.Range("J3").Value = "ee11e692-79d8-43d4-89c4-7343d3a4f709.5379b3312b0663928594db658878cb600bd5b604.zh-cn.xlf"
$zhcn.Range("K3").Value = $zhcnHandbackTime

# ---------------------------------------------------------------------
# 4. de-de sheet: same, with its own handback timestamp.
# ---------------------------------------------------------------------
$dedeHandbackTime = "2016-08-29 14:33:27"

$dede.Range("I2").Value = "a0e94d00-7c96-44aa-b06a-c48ee498b831.md"
$dede.Range("I2").Style = "Hyperlink"
$dede.Range("J2").Value = "a0e94d00-7c96-44aa-b06a-c48ee498b831.b1e5499c5bead7fe5de11bda301bce57ce7596d1.de-de.xlf"
$dede.Range("K2").Value = $dedeHandbackTime

$dede.Range("I3").Value = "ee11e692-79d8-43d4-89c4-7343d3a4f709.md"
$dede.Range("I3").Style = "Hyperlink"
$dede.Range("J3").Value = "ee11e692-79d8-43d4-89c4-7343d3a4f709.5379b3312b0663928594db658878cb600bd5b604.de-de.xlf"
$dede.Range("K3").Value = $dedeHandbackTime

# ---------------------------------------------------------------------
# 5. Hyperlinks -- add a "Latest Target File" hyperlink (column I) next to
#    the existing "Source File Name" hyperlink (column A) on each row, for
#    both sheets. Rebuild the full hyperlink list (delete + re-add, in the
#    desired A2, I2, A3, I3 order) so the new links interleave with the
#    existing ones the same way the handback report generator does.
# ---------------------------------------------------------------------
$mdTarget1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9462026a987ee4637289d88ecb367f90bff2349c/e2e/a0e94d00-7c96-44aa-b06a-c48ee498b831.md"
$mdTarget2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9462026a987ee4637289d88ecb367f90bff2349c/e2e/ee11e692-79d8-43d4-89c4-7343d3a4f709.md"

$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), $mdTarget1, "", "", "a0e94d00-7c96-44aa-b06a-c48ee498b831.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdTarget1, "", "", "a0e94d00-7c96-44aa-b06a-c48ee498b831.md")
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), $mdTarget2, "", "", "ee11e692-79d8-43d4-89c4-7343d3a4f709.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $mdTarget2, "", "", "ee11e692-79d8-43d4-89c4-7343d3a4f709.md")

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), $mdTarget1, "", "", "a0e94d00-7c96-44aa-b06a-c48ee498b831.md")
$dede.Hyperlinks.Add($dede.Range("I2"), $mdTarget1, "", "", "a0e94d00-7c96-44aa-b06a-c48ee498b831.md")
$dede.Hyperlinks.Add($dede.Range("A3"), $mdTarget2, "", "", "ee11e692-79d8-43d4-89c4-7343d3a4f709.md")
$dede.Hyperlinks.Add($dede.Range("I3"), $mdTarget2, "", "", "ee11e692-79d8-43d4-89c4-7343d3a4f709.md")
